$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Text updates (existing empty/occupied cells, styles unaffected) ---
$ws.Range("D13").Value = "Send Bias Compensated Rate Gyros"
$ws.Range("D14").Value = "Send Motor Commands 1 thru LASTMOTOR"
$ws.Range("D15").Value = "Send Motor Axis Commands"
$ws.Range("B16").Value = "Accel calibration values"
$ws.Range("D16").Value = "Send calibration values"
$ws.Range("D17").Value = "Send raw accel values"
$ws.Range("B25").Value = "Stop sending messages"
$ws.Range("D27").Value = "Send Software Configuration"

# --- Rows 28/29 take on the content that used to live in rows 30/31 ---
# C28 becomes the numeric 6 that used to be in C30 (copy number formatting first)
$ws.Range("C30").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 6

# C29 becomes the literal "=" that used to be in C31 (copy quote-prefixed formatting first)
$ws.Range("C31").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = "'="

# D28/D29 take on the text that used to live in D30/D31
$ws.Range("D28").Value = "Send Individual Motor Commands"
$ws.Range("D29").Value = "Free Form Debug"

# --- Clear out the now-vacated C30:D31 block (match C32's plain format for C31) ---
$ws.Range("C32").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C30:D31").ClearContents()

# --- Remove the trailing blank row 39 and trim C38:D38 entirely ---
$ws.Rows.Item(39).Delete()
$ws.Range("C38:D38").ClearFormats()
$ws.Range("C38:D38").ClearContents()

# --- Update selection to match the saved view state ---
$ws.Range("D20:D22").Select()
